# Adds logic to calculate a "word score" column (L) and updates the
# header row of app package ids as well as the word rows (A2:A6),
# reordering existing data and appending the new "total score" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header / package id labels ---
$ws.Range("B1").Value = "com.holoop.rollingskyball"
$ws.Range("C1").Value = "com.movisoft.rollingball"
$ws.Range("D1").Value = "com.prishka.rollsky"
$ws.Range("E1").Value = "com.besaviorgames.musicroll"
$ws.Range("F1").Value = "com.balldorsoli.rollingsky"
$ws.Range("G1").Value = "com.atesgamestudio.rollance"
$ws.Range("H1").Value = "com.dino.rolling.skyball.balance"
$ws.Range("I1").Value = "com.atreus.ballsnropes"
$ws.Range("J1").Value = "com.pronetis.gyrosphere"
$ws.Range("K1").Value = "com.doodoo.hop.ball.magic.tiles.edm.dancing.color.hooper.music.game"
$ws.Range("L1").Value = "total score"

# copy the style used by the rest of the header row onto the new L1 cell
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null

# --- Row 2 ("볼") ---
$ws.Range("A2").Value = "볼"
$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = 0.3222024913224023
$ws.Range("D2").Value = 0.2
$ws.Range("E2").Value = 0.2
$ws.Range("F2").Value = 0.2
$ws.Range("G2").Value = 0.2289440478997757
$ws.Range("H2").Value = 0.2
$ws.Range("I2").Value = 0.2
$ws.Range("J2").Value = 0.2
$ws.Range("K2").Value = 0.2
$ws.Range("L2").Value = 6.487164991287798

# --- Row 3 ("롤링") ---
$ws.Range("A3").Value = "롤링"
$ws.Range("B3").Value = 0.2
$ws.Range("C3").Value = 0.3222024913224023
$ws.Range("D3").Value = 0.2
$ws.Range("E3").Value = 0.2
$ws.Range("F3").Value = 0.2
$ws.Range("G3").Value = 0.0842238084008974
$ws.Range("H3").Value = 0.2
$ws.Range("I3").Value = 0.2
$ws.Range("J3").Value = 0.2
$ws.Range("K3").Value = 0.2
$ws.Range("L3").Value = 6.079612536257529

# --- Row 4 ("롤") ---
$ws.Range("A4").Value = "롤"
$ws.Range("B4").Value = 0.2
$ws.Range("C4").Value = 0.1185316724517319
$ws.Range("D4").Value = 0.2
$ws.Range("E4").Value = 0.2
$ws.Range("F4").Value = 0.2
$ws.Range("G4").Value = 0.2289440478997757
$ws.Range("H4").Value = 0.2
$ws.Range("I4").Value = 0.2
$ws.Range("J4").Value = 0.2
$ws.Range("K4").Value = 0.2
$ws.Range("L4").Value = 5.814007761214677

# --- Row 5 ("랜스") ---
$ws.Range("A5").Value = "랜스"
$ws.Range("B5").Value = 0.2
$ws.Range("C5").Value = 0.1185316724517319
$ws.Range("D5").Value = 0.2
$ws.Range("E5").Value = 0.2
$ws.Range("F5").Value = 0.2
$ws.Range("G5").Value = 0.2289440478997757
$ws.Range("H5").Value = 0.2
$ws.Range("I5").Value = 0.2
$ws.Range("J5").Value = 0.2
$ws.Range("K5").Value = 0.2
$ws.Range("L5").Value = 5.814007761214677

# --- Row 6 ("어드벤처") ---
$ws.Range("A6").Value = "어드벤처"
$ws.Range("B6").Value = 0.2
$ws.Range("C6").Value = 0.1185316724517319
$ws.Range("D6").Value = 0.2
$ws.Range("E6").Value = 0.2
$ws.Range("F6").Value = 0.2
$ws.Range("G6").Value = 0.2289440478997757
$ws.Range("H6").Value = 0.2
$ws.Range("I6").Value = 0.2
$ws.Range("J6").Value = 0.2
$ws.Range("K6").Value = 0.2
$ws.Range("L6").Value = 5.814007761214677
